# fix bug on message send
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Remove the resolved "Not Clearing Dialog / NewChat" issue row (old row 3).
# This shifts all following rows up by one.
$ws.Rows.Item(3).Delete()

# --- Add a new "Fixed" column (E) ------------------------------------------
# Mark the still-open issues that have now been resolved (plain centered
# style - create this style once on E3, then copy it to the other "Fixed"
# cells so only a single new style entry is produced).
$ws.Range("E3").Value = "X"
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "x"

# Row 7 (previously row 8): fix typo and mark fixed.
$ws.Range("A7").Value = "Testing needs its own environment"
$ws.Range("E3").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "x"

# Header cell, styled like the other header cells (bold white-on-grey) plus
# centered alignment.
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Fixed"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108

# --- New issue rows ----------------------------------------------------------
$ws.Range("A8").Value = "Not Switching ContextDialog "
$ws.Range("B8").Value = "Frontend"
$ws.Range("C8").Value = "DeleteChat"
$ws.Range("D8").Value = "When a Chat History is deleted the message Dialog should . Should switch to the Next Active Chat. If all chat are deleted the Message Dialog should clear"
$ws.Range("D8").WrapText = $true
$ws.Range("E3").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "X"
$ws.Rows.Item(8).RowHeight = 51

$ws.Range("A9").Value = "Not Switching ContextDialog "
$ws.Range("B9").Value = "Frontend"
$ws.Range("C9").Value = "NewChat"
$ws.Range("D9").Value = "When a new chat is created and there is an active dialog the dialof should switch to the new chat"
$ws.Range("D9").WrapText = $true
$ws.Range("E3").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "x"
$ws.Rows.Item(9).RowHeight = 34

$ws.Range("A10").Value = "Exact Message Sent Problem"
$ws.Range("B10").Value = "Frontend"
$ws.Range("C10").Value = "ChatPrompt"
$ws.Range("D10").Value = 'When the same "Exact" message is sent twice the second message is not sent'
$ws.Range("D10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 34

$ws.Range("A11").Value = "sendMessage on dup message"
$ws.Range("B11").Value = "Frontend"
$ws.Range("C11").Value = "ChatPrompt"
$ws.Range("D11").Value = "If the same message is sent twice the send of the message will not trigger a call to the AI"
$ws.Range("D11").WrapText = $true
$ws.Range("E3").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "x"
$ws.Rows.Item(11).RowHeight = 34

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 31.666666666666668
$ws.Columns.Item(5).ColumnWidth = 9.998697916666666

# --- Selection / active cell --------------------------------------------------
$ws.Range("E12").Select()
